$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new header row at the top, shifting all existing data down by one row.
$ws.Rows.Item(1).Insert()

# New header row.
$ws.Range("A1").Value = "ID"
$ws.Range("B1").Value = "DAT1_10_10"

# Fill in the 4 previously-missing DAT1 3'UTR genotypes (rows shifted down by 1
# after the header-row insert above). Force text formatting first so these
# land as shared-string cells like the rest of column B, not numbers.
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "2"
$ws.Range("B48").Style = "Normal"

$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "1"
$ws.Range("B52").Style = "Normal"

$ws.Range("B76").NumberFormat = "@"
$ws.Range("B76").Value = "1"
$ws.Range("B76").Style = "Normal"

$ws.Range("B78").NumberFormat = "@"
$ws.Range("B78").Value = "1"
$ws.Range("B78").Style = "Normal"

# Leave the selection on the new header cell (matches a plain single-cell
# selection rather than the stale whole-row selection carried over from
# the original row 1).
$ws.Range("A1").Select() | Out-Null
